# Apply the "Add new skills to database directly" edit:
#  1. Fix a typo in the "Duplicates" sheet: "Microsft Publisher" -> "Microsoft Publisher"
#  2. Re-sort Table1 on the "Redundant" sheet (column A, ascending) so the two
#     skills that had been appended at the bottom ("Depression", "Profiler")
#     land in their correct alphabetical position.

$wb = $excel.ActiveWorkbook

# --- 1. Work on the "Redundant" sheet / Table1 first -----------------------
$wsRedundant = $wb.Worksheets.Item("Redundant")
$table1 = $wsRedundant.ListObjects.Item("Table1")

$table1.Sort.SortFields.Clear()
$table1.Sort.SortFields.Add($wsRedundant.Range("A2:A43"))
$table1.Sort.Header = 1
$table1.Sort.Apply()

# Leave the selection where the user left it after sorting.
$wsRedundant.Range("A13").Select()

# --- 2. Fix the typo on the "Duplicates" sheet ------------------------------
$wsDuplicates = $wb.Worksheets.Item("Duplicates")
$wsDuplicates.Range("B65").Value = "Microsoft Publisher"

# Leave "Duplicates" as the active sheet/selection, matching the saved file.
$wsDuplicates.Range("B66").Select()
